# Updates the cryptos price/volume table to the latest scraped values.
# Leading apostrophes (doubled for PowerShell single-quote escaping) force
# Excel to store numeric-looking prices as text, matching the source data
# (preserves exact formatting like trailing zeros, e.g. "0.0690").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.111.71'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").Value = '3.256.50'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '''582.74'
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").Value = '''185.08'
$ws.Range("E6").Value = '  +1.13%  '

$ws.Range("E8").Value = '  -0.81%  '

$ws.Range("E9").Value = '  -1.79%  '

$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").Value = '3.825.06'
$ws.Range("E12").Value = '  -0.66%  '

$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").Value = '''28.15'
$ws.Range("E14").Value = '  -1.93%  '

$ws.Range("D15").Value = '68.230.77'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("D17").Value = '3.261.85'
$ws.Range("E17").Value = '  -0.77%  '

$ws.Range("D18").Value = '''5.85'
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").Value = '''393.13'
$ws.Range("E20").Value = '  +4.18%  '

$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D23").Value = '''71.34'
$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("D24").Value = '''0.519'
$ws.Range("E24").Value = '  +0.86%  '

$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("E26").Value = '  +4.24%  '

$ws.Range("D27").Value = '''9.78'
$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("E29").Value = '  -0.32%  '

$ws.Range("E30").Value = '  -0.37%  '

$ws.Range("D31").Value = '''22.91'
$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("D32").Value = '''7.16'
$ws.Range("E32").Value = '  +2.66%  '

$ws.Range("E33").Value = '  +0.10%  '

$ws.Range("D35").Value = '''1.50'
$ws.Range("E35").Value = '  -2.53%  '

$ws.Range("D36").Value = '''162.38'
$ws.Range("E36").Value = '  +0.41%  '

$ws.Range("D37").Value = '''1.95'
$ws.Range("E37").Value = '  +5.46%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '''26.93'
$ws.Range("E38").Value = '  -0.44%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").Value = '''0.823'
$ws.Range("E39").Value = '  -3.74%  '

$ws.Range("D40").Value = '''4.59'
$ws.Range("E40").Value = '  -1.05%  '

$ws.Range("E41").Value = '  -3.48%  '

$ws.Range("E42").Value = '  -6.13%  '

$ws.Range("D43").Value = '''0.0690'
$ws.Range("E43").Value = '  +1.22%  '

$ws.Range("D44").Value = '2.653.45'
$ws.Range("E44").Value = '  -0.18%  '

$ws.Range("D45").Value = '''25.41'
$ws.Range("E45").Value = '  -1.66%  '

$ws.Range("D46").Value = '''41.09'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("D47").Value = '''340.11'
$ws.Range("E47").Value = '  -3.12%  '

$ws.Range("E48").Value = '  -0.68%  '

$ws.Range("D50").Value = '''31.58'
$ws.Range("E50").Value = '  +1.41%  '

$ws.Range("D51").Value = '''0.991'
